# Daily attendance processing - 2025-12-27 07:53:49
# Swap the order of the two comma-separated "Recorded By" entries in
# column G for every row where the value starts with "dnasr281@gmail.com, ".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value()

    if ($val -ne $null -and $val -is [string] -and $val.Contains(", ")) {
        $parts = $val -split ", ", 2
        if ($parts.Count -eq 2 -and $parts[0] -eq "dnasr281@gmail.com") {
            $cell.Value = $parts[1] + ", " + $parts[0]
        }
    }
}
